$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Straightforward stat refreshes (same country, same row, just new numbers)
# ---------------------------------------------------------------------------

# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value = 5999676
$ws.Cells.Item(4,3).Value = 43948
$ws.Cells.Item(4,4).Value = 3296352
$ws.Cells.Item(4,5).Value = 2519683
$ws.Cells.Item(4,7).Value = 1277
$ws.Cells.Item(4,8).Value = 183641

# Row 9: Peru
$ws.Cells.Item(9,2).Value = 613378
$ws.Cells.Item(9,3).Value = 5996
$ws.Cells.Item(9,4).Value = 421877
$ws.Cells.Item(9,5).Value = 163377
$ws.Cells.Item(9,7).Value = 123
$ws.Cells.Item(9,8).Value = 28124

# Row 115: Suazilandia
$ws.Cells.Item(115,2).Value = 4387
$ws.Cells.Item(115,3).Value = 60
$ws.Cells.Item(115,4).Value = 3078
$ws.Cells.Item(115,5).Value = 1221
$ws.Cells.Item(115,7).Value = 2
$ws.Cells.Item(115,8).Value = 88

# Row 165: Santo Tome y Principe
$ws.Cells.Item(165,4).Value = 842
$ws.Cells.Item(165,5).Value = 35

# Row 190: Barbados
$ws.Cells.Item(190,2).Value = 165
$ws.Cells.Item(190,3).Value = 1
$ws.Cells.Item(190,4).Value = 139
$ws.Cells.Item(190,5).Value = 19

# Row 193: Monaco
$ws.Cells.Item(193,2).Value = 122
$ws.Cells.Item(193,3).Value = 1
$ws.Cells.Item(193,5).Value = 36

# ---------------------------------------------------------------------------
# 2) Aruba moves up the ranking: it used to sit right after "Nueva Zelanda"
#    (row 147); now it belongs right after "Bahamas" (row 142), pushing
#    "Jordania", "Malta", "Jamaica" and "Nueva Zelanda" down by one row.
# ---------------------------------------------------------------------------

$ws.Rows.Item(147).Delete()
$ws.Rows.Item(143).Insert()
$ws.Cells.Item(143,1).Value = "Aruba"
$ws.Cells.Item(143,2).Value = 1760
$ws.Cells.Item(143,3).Value = 90
$ws.Cells.Item(143,4).Value = 587
$ws.Cells.Item(143,5).Value = 1165
$ws.Cells.Item(143,6).Value = 0
$ws.Cells.Item(143,7).Value = 0
$ws.Cells.Item(143,8).Value = 8

# ---------------------------------------------------------------------------
# 3) Martinica moves up the ranking: it used to sit right after "Taiwan"
#    (row 172); now it belongs right after "San Marino" (row 168), pushing
#    "Birmania", "Tanzania" and "Taiwan" down by one row.
# ---------------------------------------------------------------------------

$ws.Rows.Item(172).Delete()
$ws.Rows.Item(169).Insert()
$ws.Cells.Item(169,1).Value = "Martinica"
$ws.Cells.Item(169,2).Value = 596
$ws.Cells.Item(169,3).Value = 132
$ws.Cells.Item(169,4).Value = 98
$ws.Cells.Item(169,5).Value = 482
$ws.Cells.Item(169,6).Value = 0
$ws.Cells.Item(169,7).Value = 0
$ws.Cells.Item(169,8).Value = 16

# ---------------------------------------------------------------------------
# 4) Refresh the "last updated" timestamp caption in A1
# ---------------------------------------------------------------------------

$ws.Cells.Item(1,1).Value = "Datos actualizados a 27 de Agosto de 2020 a las 03:04"
